# Insert a new data row at row 103 (pushing the existing rows 103..199 down to
# 104..200) and populate it with the new "Fruta / hortaliza, semanal" record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(103).Insert()

$ws.Cells.Item(103, 1).Value = 4
$ws.Cells.Item(103, 2).Value = 'Feria Lagunitas de Puerto Montt'
$ws.Cells.Item(103, 3).Value = 'Los Lagos'
$ws.Cells.Item(103, 4).Value = 44566
$ws.Cells.Item(103, 5).Value = 10
$ws.Cells.Item(103, 6).Value = 'Fruta'
$ws.Cells.Item(103, 7).Value = 100102
$ws.Cells.Item(103, 8).Value = 'Cítricos'
$ws.Cells.Item(103, 9).Value = 100102006
$ws.Cells.Item(103, 10).Value = 'Pomelo'
$ws.Cells.Item(103, 11).Value = 'Start Ruby'
$ws.Cells.Item(103, 12).Value = 'Primera'
$ws.Cells.Item(103, 13).Value = 30
$ws.Cells.Item(103, 14).Value = 11000
$ws.Cells.Item(103, 15).Value = 12000
$ws.Cells.Item(103, 16).Value = 11500
$ws.Cells.Item(103, 17).Value = '$/caja 14 kilos empedrada'
$ws.Cells.Item(103, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(103, 19).Value = 821
$ws.Cells.Item(103, 20).Value = 14
